$wb = $excel.ActiveWorkbook
Write-Output $wb.Path
Write-Output $wb.FullName
try { $wb.Path = "D:\Daten\t1057\Desktop\Local Working Directory\07_Projekte\Lp Corruptions\Github Code and Results\results\"; Write-Output "set ok" } catch { Write-Output "ERR: $_" }
